$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$participants_q = 'MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE diag.primary_diagnosis in [''Anaplastic medulloblastoma'']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'''') as `Participant ID`,
coalesce(s.study_name, '''') as `Study Name`,
coalesce(s.phs_accession,'''') as `Accession`,
coalesce(p.gender,'''') as `Gender`,
coalesce(apoc.text.join(samp, '',''), '''') as `Samples`
ORDER BY p.participant_id
LIMIT 100'
$files_q = 'MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE diag.primary_diagnosis in [''Anaplastic medulloblastoma'']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '''') as `File Name`,
    coalesce(s.study_name,'''') as `Study Name`,
    coalesce(s.phs_accession,'''') as `Accession`,
    coalesce(p.participant_id, '''') as `Participant ID`,
    coalesce(samp.sample_id, '''') as `Sample ID`,
    coalesce(f.file_type, '''') as `File Type`
ORDER BY f.file_name
LIMIT 100'
$samples_q = 'MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE diag.primary_diagnosis in [''Anaplastic medulloblastoma'']
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '''') as `Sample ID`,
    coalesce(p.participant_id,'''') as `Participant ID`,
    coalesce(s.study_name, '''') as `Study Name`,
    coalesce(s.phs_accession,'''') as `Accession`,
    coalesce(samp.sample_tumor_status,'''') as `Tumor`,
    coalesce(samp.sample_type,'''') as `Analyte Type`
ORDER BY samp.sample_id 
LIMIT 100'
$statquery_q = 'CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
   WHERE diag.primary_diagnosis in [''Anaplastic medulloblastoma'']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE diag.primary_diagnosis in [''Anaplastic medulloblastoma'']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE diag.primary_diagnosis in [''Anaplastic medulloblastoma'']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`'

$ws.Range("B2").Value = $participants_q
$ws.Range("C2").Value = $statquery_q
$ws.Range("B3").Value = $samples_q
$ws.Range("C3").Value = $statquery_q
$ws.Range("B4").Value = $files_q
$ws.Range("C4").Value = $statquery_q

$ws.Rows.Item(4).RowHeight = 409.5

$wb.Save()
